$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 58
$ws.Range("I2").Value = 130
$ws.Range("J2").Value = 626
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 155
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 116
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 3
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 54
$ws.Range("T2").Value = 116
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 969
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 968
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 4
